$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new job posting row (row 4) with Job_Id = JD_003
$ws.Range("A4").Value = "JD_003"
$ws.Range("B4").Value = "Cyber Security Engineer"
$ws.Range("C4").Value = "Testing"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 4
